# MEDICIONES BERLANGA-ANDALUZ: rellenar la columna de fechas (L) para cada
# capataz, eliminando los datos de prueba duplicados ("PRUEBAS", "SADÑLNSA",
# "GAM", "DOS", "SALNSA", "SALSA") y sustituyéndolos por los días del mes
# (15 a 28 de abril, en letra) que faltaban.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L5").Value  = "QUINCE"
$ws.Range("L6").Value  = "DIECISEIS"
$ws.Range("L7").Value  = "DIECISIETE"
$ws.Range("L8").Value  = "DIECIOCHO"
$ws.Range("L9").Value  = "DIECINUEVE"
$ws.Range("L10").Value = "VEINTE"
$ws.Range("L11").Value = "VEINTIUNO"
$ws.Range("L12").Value = "VEINTIDOS"
$ws.Range("L13").Value = "VIENTITRES"
$ws.Range("L14").Value = "VEINTICUATRO"
$ws.Range("L15").Value = "VEINTICINCO"
$ws.Range("L16").Value = "VEINTISEIS"
$ws.Range("L17").Value = "VEINTISIETE"
$ws.Range("L18").Value = "VEINTIOCHO"

# La selección activa queda en L18, como en el libro final.
$ws.Range("L18").Select() | Out-Null
